$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 50, pushing the existing rows 50-62 down to 52-64.
$ws.Rows("50:51").Insert()

# New row 50: "Región de Ñuble", Primera, 60, 5000/5000/5000, 2500
$ws.Cells.Item(50, 1).Value = 7
$ws.Cells.Item(50, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(50, 3).Value = "Ñuble"
$ws.Cells.Item(50, 4).Value = 45275
$ws.Cells.Item(50, 5).Value = 16
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100101
$ws.Cells.Item(50, 8).Value = "Berries"
$ws.Cells.Item(50, 9).Value = 100101001
$ws.Cells.Item(50, 10).Value = "Arándano (blue)"
$ws.Cells.Item(50, 11).Value = "Sin especificar"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 60
$ws.Cells.Item(50, 14).Value = 5000
$ws.Cells.Item(50, 15).Value = 5000
$ws.Cells.Item(50, 16).Value = 5000
$ws.Cells.Item(50, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(50, 18).Value = "Región de Ñuble"
$ws.Cells.Item(50, 19).Value = 2500
$ws.Cells.Item(50, 20).Value = 2

# New row 51: "Región de Ñuble", Segunda, 60, 4000/4000/4000, 2000
$ws.Cells.Item(51, 1).Value = 7
$ws.Cells.Item(51, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(51, 3).Value = "Ñuble"
$ws.Cells.Item(51, 4).Value = 45275
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(51, 6).Value = "Fruta"
$ws.Cells.Item(51, 7).Value = 100101
$ws.Cells.Item(51, 8).Value = "Berries"
$ws.Cells.Item(51, 9).Value = 100101001
$ws.Cells.Item(51, 10).Value = "Arándano (blue)"
$ws.Cells.Item(51, 11).Value = "Sin especificar"
$ws.Cells.Item(51, 12).Value = "Segunda"
$ws.Cells.Item(51, 13).Value = 60
$ws.Cells.Item(51, 14).Value = 4000
$ws.Cells.Item(51, 15).Value = 4000
$ws.Cells.Item(51, 16).Value = 4000
$ws.Cells.Item(51, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(51, 18).Value = "Región de Ñuble"
$ws.Cells.Item(51, 19).Value = 2000
$ws.Cells.Item(51, 20).Value = 2
